$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: 91489339 -> 75185465
$ws.Range("A2").Value = 75185465

# S2: 25 -> 10
$ws.Range("S2").Value = 10

# AF2: new empty (text) cell, mirrors the blank cells already present at
# I2/AT2 (empty inline string). A bare "" assignment is a no-op for the
# COM model (it leaves the cell absent), so force an empty *text* entry
# with a leading apostrophe, matching how Excel persists a deliberately
# blank text cell.
$ws.Range("AF2").Value = "'"

# AW2 (Rapportör) and AX2 (Observatörer): both become "Sofie Jonsson"
$ws.Range("AW2").Value = "Sofie Jonsson"
$ws.Range("AX2").Value = "Sofie Jonsson"
